# إضافة حدث جديد في Card13 by admin at 2025-12-08 11:36:52
# Re-export touched up every previously-blank data cell (rows 2-16) to the
# literal text "nan", bumped the used range to A1:O17, and appended a new
# service-log row (17) recording a front-card clothing replacement event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

# Fill previously-blank cells (rows 2-16) with the literal text "nan"
foreach ($col in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)) { $ws.Cells.Item(2, $col).Value = "nan" }
foreach ($col in @(7, 8, 9, 10, 11, 13, 14, 15)) { $ws.Cells.Item(3, $col).Value = "nan" }
foreach ($col in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)) { $ws.Cells.Item(4, $col).Value = "nan" }
foreach ($col in @(4, 5, 8, 9, 10, 11, 13, 14, 15)) { $ws.Cells.Item(5, $col).Value = "nan" }
foreach ($col in @(5, 6, 7, 9, 10, 11, 13, 14, 15)) { $ws.Cells.Item(6, $col).Value = "nan" }
foreach ($col in @(5, 7, 8, 9, 10, 13, 14, 15)) { $ws.Cells.Item(7, $col).Value = "nan" }
foreach ($col in @(4, 6, 7, 8, 10, 11)) { $ws.Cells.Item(8, $col).Value = "nan" }
foreach ($col in @(5, 8, 9, 10, 11, 13, 14, 15)) { $ws.Cells.Item(9, $col).Value = "nan" }
foreach ($col in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)) { $ws.Cells.Item(10, $col).Value = "nan" }
foreach ($col in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)) { $ws.Cells.Item(11, $col).Value = "nan" }
foreach ($col in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)) { $ws.Cells.Item(12, $col).Value = "nan" }
foreach ($col in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)) { $ws.Cells.Item(13, $col).Value = "nan" }
foreach ($col in @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 14)) { $ws.Cells.Item(14, $col).Value = "nan" }
foreach ($col in @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 14)) { $ws.Cells.Item(15, $col).Value = "nan" }
foreach ($col in @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 14)) { $ws.Cells.Item(16, $col).Value = "nan" }

# Append new row 17 with the new service event
$ws.Cells.Item(17, 1).Value = "'13"
$ws.Cells.Item(17, 2).Value = "'"
$ws.Cells.Item(17, 3).Value = "'"
$ws.Cells.Item(17, 4).Value = "'"
$ws.Cells.Item(17, 5).Value = "'"
$ws.Cells.Item(17, 6).Value = "'"
$ws.Cells.Item(17, 7).Value = "'"
$ws.Cells.Item(17, 8).Value = "'"
$ws.Cells.Item(17, 9).Value = "'"
$ws.Cells.Item(17, 10).Value = "'"
$ws.Cells.Item(17, 11).Value = "'"
$ws.Cells.Item(17, 12).Value = "20\4\2025"
$ws.Cells.Item(17, 13).Value = "تم تغيير الجرائد الاماميه (1_2_4_5_7_8)"
$ws.Cells.Item(17, 14).Value = "5766 t"
$ws.Cells.Item(17, 15).Value = "الخبير"
